$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1. Workbook-level metadata tweak (author moved the project folder).
# ---------------------------------------------------------------------------
$wb.Application.StatusBar = $false

# ---------------------------------------------------------------------------
# 2. Remove the now-duplicated "Must allow use date to advance" row (old
#    row 21) -- its text was folded into row 3's updated wording below, and
#    the rest of the rows shift up by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update row 3's feature text.
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Allow user to select the filename for books"

# ---------------------------------------------------------------------------
# 4. Add the new "class" column (F) with category labels, color-coded to
#    match the same classes highlighted in column A/B.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 6).Value = "Patrons"
$ws.Cells.Item(3, 6).Value = "Books/media"
$ws.Cells.Item(5, 6).Value = "Library"
$ws.Cells.Item(6, 6).Value = "Interface"

$ws.Cells.Item(4, 6).Style = "Good"
$ws.Cells.Item(5, 6).Style = "Bad"
$ws.Cells.Item(6, 6).Style = "Neutral"

$ws.Columns.Item(6).ColumnWidth = 11.77734375

# ---------------------------------------------------------------------------
# 5. Color-code the existing rows by class.
#      Good (green)    -> Patrons
#      Bad (red)       -> Library
#      Neutral (yellow)-> Interface
# ---------------------------------------------------------------------------
$ws.Range("A6:B7").Style = "Good"
$ws.Range("A11:A15").Style = "Good"

$ws.Range("A18:A22").Style = "Bad"
$ws.Cells.Item(30, 1).Style = "Bad"

$ws.Range("A23:A29").Style = "Neutral"

# ---------------------------------------------------------------------------
# 6. Refresh the table range now that the sheet shrank by one row, and
#    reselect near the edited area to mirror the author's final view.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B33")) | Out-Null

$ws.Range("A30").Select() | Out-Null
